$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) holds numeric-looking text (e.g. "68.60", "4.10")
# that must stay literal text -- Excel would otherwise auto-convert it to a
# number and silently drop meaningful trailing zeros. Force text format first.

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "36.388.14"
$ws.Range("E2").Value = "  -0.15%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.932.05"
$ws.Range("E3").Value = "  -2.69%  "

# Row 4
$ws.Range("E4").Value = "  +0.07%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.69"
$ws.Range("E5").Value = "  -2.04%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.606"
$ws.Range("E6").Value = "  -3.84%  "

# Row 7
$ws.Range("E7").Value = "  -0.02%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "56.07"
$ws.Range("E8").Value = "  -6.79%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.356"
$ws.Range("E9").Value = "  -5.39%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0834"
$ws.Range("E10").Value = "  +1.16%  "

# Row 11
$ws.Range("E11").Value = "  -1.07%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.219.32"
$ws.Range("E12").Value = "  -2.42%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.797"
$ws.Range("E13").Value = "  -8.36%  "

# Row 14
$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "13.27"
$ws.Range("E14").Value = "  -5.86%  "

# Row 15
$ws.Range("B15").Value = "Avalanche"
$ws.Range("C15").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "20.77"
$ws.Range("E15").Value = "  -12.80%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.09"
$ws.Range("E16").Value = "  -7.22%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.943.63"
$ws.Range("E17").Value = "  -1.89%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "36.281.00"
$ws.Range("E18").Value = "  -0.05%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0₃0857"
$ws.Range("E19").Value = "  -2.21%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "68.60"
$ws.Range("E20").Value = "  -2.67%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "225.94"
$ws.Range("E21").Value = "  -3.56%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.91"
$ws.Range("E22").Value = "  -8.25%  "

# Row 23
$ws.Range("E23").Value = "  -0.22%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.32"
$ws.Range("E24").Value = "  -10.36%  "

# Row 25
$ws.Range("E25").Value = "  -3.06%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.18"
$ws.Range("E26").Value = "  -8.46%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "160.23"
$ws.Range("E27").Value = "  -1.38%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.131"
$ws.Range("E28").Value = "  -0.15%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.04"
$ws.Range("E29").Value = "  -4.33%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.116"
$ws.Range("E30").Value = "  -3.25%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.09"
$ws.Range("E31").Value = "  -7.87%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.51"
$ws.Range("E32").Value = "  -8.50%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0621"
$ws.Range("E33").Value = "  -4.87%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.10"
$ws.Range("E34").Value = "  -7.59%  "

# Row 35
$ws.Range("E35").Value = "  +0.09%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.01"
$ws.Range("E36").Value = "  -3.45%  "

# Row 37
$ws.Range("E37").Value = "  -0.03%  "

# Row 38
$ws.Range("E38").Value = "  -7.26%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.94"
$ws.Range("E39").Value = "  -1.76%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0964"
$ws.Range("E40").Value = "  -1.32%  "

# Row 41
$ws.Range("E41").Value = "  -1.08%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0207"
$ws.Range("E42").Value = "  -3.43%  "

# Row 43
$ws.Range("E43").Value = "  -8.51%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "15.37"
$ws.Range("E44").Value = "  -5.89%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.325.79"
$ws.Range("E45").Value = "  -3.48%  "

# Row 46
$ws.Range("E46").Value = "  -8.13%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "84.99"
$ws.Range("E47").Value = "  -8.35%  "

# Row 48
$ws.Range("E48").Value = "  -6.64%  "

# Row 49
$ws.Range("E49").Value = "  -0.31%  "

# Row 50
$ws.Range("B50").Value = "RocketPoolETH"
$ws.Range("C50").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.109.79"
$ws.Range("E50").Value = "  -2.45%  "

# Row 51
$ws.Range("B51").Value = "MultiversX"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "43.24"
$ws.Range("E51").Value = "  -5.40%  "
